$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D column) values that only change price, no row reordering ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "246.63"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.46"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.255"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05684"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.418"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.308"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8062"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8601"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1413"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07396"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03038"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006812"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1068"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008452"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005584"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4502"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1600"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01010"

# --- Rows 14-26 got reordered (row 26 ProBitToken moved to row 14, others shifted down by one) ---
# and each row also received an updated Price (D) value.
$ws.Range("B14").Value = "ProBitToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1282"
$ws.Range("E14").Value = "13ProBitTokenPROB"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09385"
$ws.Range("E15").Value = "14BitMartTokenBMX"

$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.867"
$ws.Range("E16").Value = "15MCDexMCB"

$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001583"
$ws.Range("E17").Value = "16BitForexTokenBF"

$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04784"
$ws.Range("E18").Value = "17CoinExTokenCET"

$ws.Range("B19").Value = "One"
$ws.Range("C19").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0005853"
$ws.Range("E19").Value = "18OneONE"

$ws.Range("B20").Value = "TigerCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.006387"
$ws.Range("E20").Value = "19TigerCashTCH"

$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.005033"
$ws.Range("E21").Value = "20HotbitTokenHTB"

$ws.Range("B22").Value = "BitKan"
$ws.Range("C22").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0009973"
$ws.Range("E22").Value = "21BitKanKAN"

$ws.Range("B23").Value = "NitroEx"
$ws.Range("C23").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0001500"
$ws.Range("E23").Value = "22NitroExNTX"

$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.689"
$ws.Range("E24").Value = "23LEOLEO"

$ws.Range("B25").Value = "BTSEToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.198"
$ws.Range("E25").Value = "24BTSETokenBTSE"

$ws.Range("B26").Value = "BitpandaEcosystemToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.3279"
$ws.Range("E26").Value = "25BitpandaEcosystemTokenBEST"

